# Add two new portfolio rows for a new holder "HARSH VARDHAN" (STERTOOLS, TATAMOTORS)
# and give the AVG_PRICE column (D) a left-aligned number style for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 64 (values entered column by column, C before B, matching the
# order the new shared strings were appended to the workbook)
$ws.Cells.Item(64, 3).Value = "STERTOOLS"
$ws.Cells.Item(64, 2).Value = "HARSH VARDHAN"
$ws.Cells.Item(64, 1).Value = 62
$ws.Cells.Item(64, 4).Value = 356.45
$ws.Cells.Item(64, 5).Value = 5
$ws.Cells.Item(64, 6).Value = 3
$ws.Cells.Item(64, 7).Value = "nunna.harshavardhan2001@gmail.com"

# New row 65
$ws.Cells.Item(65, 3).Value = "TATAMOTORS"
$ws.Cells.Item(65, 2).Value = "HARSH VARDHAN"
$ws.Cells.Item(65, 1).Value = 63
$ws.Cells.Item(65, 4).Value = 873.5667
$ws.Cells.Item(65, 5).Value = 3
$ws.Cells.Item(65, 6).Value = 3
$ws.Cells.Item(65, 7).Value = "nunna.harshavardhan2001@gmail.com"

# Copy the numbering style used by column A (bold/border/centered) onto the new rows
$ws.Cells.Item(62, 1).Copy() | Out-Null
$ws.Range("A64:A65").PasteSpecial(-4122) | Out-Null

# Give the AVG_PRICE cells a left-aligned style
$ws.Range("D64:D65").HorizontalAlignment = -4131

# Update the visible window / selection to match the edited state
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("J60").Select() | Out-Null
